# Add a new weekly-report block (rows 261-270) to Sheet1, mirroring the
# structure of the immediately preceding block (rows 251-260):
#   - 1 merged "date/section" header row
#   - 1 merged "column header" row (组员/计划内容/完成情况/备注)
#   - 6 per-member rows
#   - 1 merged 2-row "总结：" (summary) block
#
# We clone formatting by copying each source cell individually (not as a
# merged range) into the matching destination cell — this preserves the
# exact existing cell style (no new style/border combinations get minted)
# the way a whole-range merged-cell copy would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcStart = 251
$dstStart = 261
$blockRows = 10

for ($r = 0; $r -lt $blockRows; $r++) {
    $srcRow = $srcStart + $r
    $dstRow = $dstStart + $r
    for ($c = 1; $c -le 4; $c++) {
        $ws.Cells.Item($srcRow, $c).Copy($ws.Cells.Item($dstRow, $c))
    }
}

# Re-create the merges for the new block (same relative layout as A251:D260):
$ws.Range("A261:D261").Merge()
$ws.Range("C263:C268").Merge()
$ws.Range("D263:D268").Merge()
$ws.Range("A269:D270").Merge()

# --- Fill in the new block's content -----------------------------------

# Section header / date row
$ws.Range("A261").Value = "日期：2018.12.3第十四周周一"

# Column header row (unchanged boilerplate, already copied, but set explicitly
# to be safe / explicit about intent)
$ws.Range("A262").Value = "组员"
$ws.Range("B262").Value = "计划内容"
$ws.Range("C262").Value = "完成情况"
$ws.Range("D262").Value = "备注"

# Per-member rows
$ws.Range("A263").Value = "陈柯赞"
$ws.Range("B263").Value = "测试报告"

$ws.Range("A264").Value = "黎安生"
$ws.Range("B264").Value = "网页端使用手册"

$ws.Range("A265").Value = "王智永"
$ws.Range("B265").Value = "APP端完善"

$ws.Range("A266").Value = "郑海文"
$ws.Range("B266").Value = "ppt"

$ws.Range("A267").Value = "赵华亮"
$ws.Range("B267").Value = "app端使用手册"

$ws.Range("A268").Value = "叶田"
$ws.Range("B268").Value = "辅助"

# Clear leftover values copied from the source block's 完成情况/备注 columns
$ws.Range("C263").Value = ""
$ws.Range("D263").Value = ""
$ws.Range("C264").Value = ""
$ws.Range("D264").Value = ""
$ws.Range("C265").Value = ""
$ws.Range("D265").Value = ""
$ws.Range("C266").Value = ""
$ws.Range("D266").Value = ""
$ws.Range("C267").Value = ""
$ws.Range("D267").Value = ""
$ws.Range("C268").Value = ""
$ws.Range("D268").Value = ""

# Summary row
$ws.Range("A269").Value = "总结："
$ws.Range("B269").Value = ""
$ws.Range("C269").Value = ""
$ws.Range("D269").Value = ""
$ws.Range("A270").Value = ""
$ws.Range("B270").Value = ""
$ws.Range("C270").Value = ""
$ws.Range("D270").Value = ""

# Move the active selection to the newly-added summary block, matching how
# the previous edit left the selection on the then-last block (A259:D260).
$ws.Range("A269:D270").Select()
